$wb = $excel.ActiveWorkbook

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3005.1428
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 3005.1428
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 3005.1428
$ws.Range("N113").Value = -9513.1428
$ws.Range("M113").ClearContents()

# ALC row 117
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# ALC row 123
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 40000
$ws.Range("J123").Value = 40000
$ws.Range("L123").Value = 40000
$ws.Range("N123").Value = -49800

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2435.2537
$ws.Range("I132").Value = 2329.242
$ws.Range("J132").Value = 3749.8
$ws.Range("K132").Value = 6987.726000000001
$ws.Range("L132").Value = 11249.4
$ws.Range("M132").Value = -4457.726000000001
$ws.Range("N132").Value = -16309.4

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1168.3898
$ws.Range("I137").Value = 956.16327
$ws.Range("J137").Value = 2208.3
$ws.Range("K137").Value = 2868.48981
$ws.Range("L137").Value = 6624.900000000001
$ws.Range("M137").Value = -318.48981
$ws.Range("N137").Value = -11724.9

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2184.5652
$ws.Range("I138").Value = 1483.9744
$ws.Range("J138").Value = 2700.0942
$ws.Range("K138").Value = 4451.9232
$ws.Range("L138").Value = 8100.2826
$ws.Range("M138").Value = 688.0767999999998
$ws.Range("N138").Value = -18380.2826

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2439.386
$ws.Range("I141").Value = 869.7917
$ws.Range("J141").Value = 10810.556
$ws.Range("K141").Value = 2609.3751
$ws.Range("L141").Value = 32431.668
$ws.Range("M141").Value = 2570.6249
$ws.Range("N141").Value = -42791.66800000001

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 750616.4
$ws.Range("I32").Value = 827211.2
$ws.Range("J32").Value = 22965.875
$ws.Range("K32").Value = 827211.2
$ws.Range("L32").Value = 22965.875
$ws.Range("M32").Value = -826924.2
$ws.Range("N32").Value = -23539.875

# ARM row 34
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1370.5333
$ws.Range("I74").Value = 1126.4615
$ws.Range("K74").Value = 1126.4615
$ws.Range("M74").Value = -252.4614999999999

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1370.5333
$ws.Range("I77").Value = 1126.4615
$ws.Range("K77").Value = 5632.307499999999
$ws.Range("M77").Value = -1264.307499999999

# ARM row 86
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H86").Value = 50028500
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 50028500
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 50028500
$ws.Range("N86").Value = -50030872
$ws.Range("M86").ClearContents()

# ARM row 89
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H89").Value = 50028500
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 50028500
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 150085500
$ws.Range("N89").Value = -150097356
$ws.Range("M89").ClearContents()

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3202
$ws.Range("I102").Value = 2752.5
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 2752.5
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -1130.5
$ws.Range("N102").Value = -8244

# ARM row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 44444
$ws.Range("J113").Value = 44444
$ws.Range("L113").Value = 44444
$ws.Range("N113").Value = -53122

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3240.4043
$ws.Range("I132").Value = 2441.9355
$ws.Range("K132").Value = 7325.806500000001
$ws.Range("M132").Value = -4795.806500000001

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2786.639
$ws.Range("I134").Value = 2319.9614
$ws.Range("K134").Value = 6959.8842
$ws.Range("M134").Value = -4424.8842

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 931.9091
$ws.Range("I58").Value = 648.4865
$ws.Range("J58").Value = 1514.5
$ws.Range("K58").Value = 648.4865
$ws.Range("L58").Value = 1514.5
$ws.Range("M58").Value = -445.4865
$ws.Range("N58").Value = -1920.5

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3473529.5
$ws.Range("I132").Value = 1151.25
$ws.Range("J132").Value = 13890664
$ws.Range("K132").Value = 3453.75
$ws.Range("L132").Value = 41671992
$ws.Range("M132").Value = -923.75
$ws.Range("N132").Value = -41677052

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3116
$ws.Range("I134").Value = 3045.5625
$ws.Range("J134").Value = 3491.6667
$ws.Range("K134").Value = 9136.6875
$ws.Range("L134").Value = 10475.0001
$ws.Range("M134").Value = -6601.6875
$ws.Range("N134").Value = -15545.0001

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 931.9091
$ws.Range("I136").Value = 648.4865
$ws.Range("J136").Value = 1514.5
$ws.Range("K136").Value = 1945.4595
$ws.Range("L136").Value = 4543.5
$ws.Range("M136").Value = 604.5405000000001
$ws.Range("N136").Value = -9643.5

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1231.575
$ws.Range("J5").Value = 1868.9131
$ws.Range("L5").Value = 5606.7393
$ws.Range("N5").Value = -5830.7393

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1231.575
$ws.Range("J135").Value = 1868.9131
$ws.Range("L135").Value = 16820.2179
$ws.Range("N135").Value = -21890.2179

# LTW row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1483.2325
$ws.Range("I68").Value = 1442.6757
$ws.Range("J68").Value = 1733.3334
$ws.Range("K68").Value = 1442.6757
$ws.Range("L68").Value = 1733.3334
$ws.Range("M68").Value = -693.6757
$ws.Range("N68").Value = -3231.3334

# LTW row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 1483.2325
$ws.Range("I71").Value = 1442.6757
$ws.Range("J71").Value = 1733.3334
$ws.Range("K71").Value = 7213.3785
$ws.Range("L71").Value = 8666.666999999999
$ws.Range("M71").Value = -3469.3785
$ws.Range("N71").Value = -16154.667

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3402667.2
$ws.Range("I136").Value = 1232.0968
$ws.Range("J136").Value = 9260695
$ws.Range("K136").Value = 3696.2904
$ws.Range("L136").Value = 27782085
$ws.Range("M136").Value = -1146.2904
$ws.Range("N136").Value = -27787185

# WVR row 117
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2424.3777
$ws.Range("I122").Value = 2147.6
$ws.Range("J122").Value = 2977.9333
$ws.Range("K122").Value = 6442.799999999999
$ws.Range("L122").Value = 8933.7999
$ws.Range("M122").Value = -3992.799999999999
$ws.Range("N122").Value = -13833.7999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2348206.8
$ws.Range("I132").Value = 704.69385
$ws.Range("J132").Value = 7576734
$ws.Range("K132").Value = 2114.08155
$ws.Range("L132").Value = 22730202
$ws.Range("M132").Value = 415.9184500000001
$ws.Range("N132").Value = -22735262
